$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp column (O) for all data rows (2 through 73)
for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-08-15 20:57:31"
}

# Row 8: Avela Strumpfhose - price drops from 3.95 to 1.95 (50% off)
$ws.Cells.Item(8, 8).NumberFormat = "@"
$ws.Cells.Item(8, 8).Value = "1.95"
$ws.Cells.Item(8, 14).Value = "Avela Strumpfhose Top Size Noir  11 - 12 50% Aktion 1.95 Schweizer Franken statt 3.95 Schweizer Franken"

# Row 69: Selenacare Menstruationsunterwäsche - price drops from 24.50 to 12.25 (50% off)
$ws.Cells.Item(69, 8).NumberFormat = "@"
$ws.Cells.Item(69, 8).Value = "12.25"
$ws.Cells.Item(69, 9).Value = "12.25/1ST"
$ws.Cells.Item(69, 11).NumberFormat = "@"
$ws.Cells.Item(69, 11).Value = "12.25"
$ws.Cells.Item(69, 14).Value = "Selenacare Menstruationsunterwäsche S 50% Aktion 12.25 Schweizer Franken statt 24.50 Schweizer Franken"
